$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.073391145793255
$ws.Range("D2").Value = 1.058725556297544
$ws.Range("E2").Value = 1.07437935542267
$ws.Range("F2").Value = 1.077832230460339
$ws.Range("I2").Value = 1.043165514517731
$ws.Range("J2").Value = 1.078305890488911
$ws.Range("K2").Value = 1.06145700513595
$ws.Range("L2").Value = 1.077068723531604
$ws.Range("M2").Value = 1.080512494309527
$ws.Range("N2").Value = 1.07983720853682
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.075658260135517
$ws.Range("D3").Value = 1.059703262075876
$ws.Range("E3").Value = 1.076244415215817
$ws.Range("F3").Value = 1.079397571625208
$ws.Range("I3").Value = 1.043449312399833
$ws.Range("J3").Value = 1.080225538171354
$ws.Range("K3").Value = 1.062248385890024
$ws.Range("L3").Value = 1.078748326518776
$ws.Range("M3").Value = 1.081893780964458
$ws.Range("N3").Value = 1.081759582339159
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.077120534065552
$ws.Range("D4").Value = 1.060333136967861
$ws.Range("E4").Value = 1.077447007264279
$ws.Range("F4").Value = 1.080406500899334
$ws.Range("I4").Value = 1.043630431142659
$ws.Range("J4").Value = 1.081462857827303
$ws.Range("K4").Value = 1.062757084080051
$ws.Range("L4").Value = 1.079830485781488
$ws.Range("M4").Value = 1.082783120617612
$ws.Range("N4").Value = 1.082998659130942
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.077734178878304
$ws.Range("D5").Value = 1.060597283125925
$ws.Range("E5").Value = 1.077951588396974
$ws.Range("F5").Value = 1.080829726771562
$ws.Range("I5").Value = 1.043705974614997
$ws.Range("J5").Value = 1.081981898294406
$ws.Range("K5").Value = 1.062970140851216
$ws.Range("L5").Value = 1.080284332702117
$ws.Range("M5").Value = 1.083155951790691
$ws.Range("N5").Value = 1.083518436695043
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.077837149195752
$ws.Range("D6").Value = 1.060641596380569
$ws.Range("E6").Value = 1.078036252491302
$ws.Range("F6").Value = 1.080900734357653
$ws.Range("I6").Value = 1.0437186236906
$ws.Range("J6").Value = 1.082068982018195
$ws.Range("K6").Value = 1.063005867364639
$ws.Range("L6").Value = 1.080360472145054
$ws.Range("M6").Value = 1.083218490866356
$ws.Range("N6").Value = 1.083605644087708
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.077128737880084
$ws.Range("D7").Value = 1.060336669056305
$ws.Range("E7").Value = 1.077453753355084
$ws.Range("F7").Value = 1.080412159686267
$ws.Range("I7").Value = 1.043631442906038
$ws.Range("J7").Value = 1.081469797672279
$ws.Range("K7").Value = 1.062759934086349
$ws.Range("L7").Value = 1.079836554370819
$ws.Range("M7").Value = 1.082788106497865
$ws.Range("N7").Value = 1.083005608831294
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.074158319283526
$ws.Range("D8").Value = 1.059056554724815
$ws.Range("E8").Value = 1.075010550273976
$ws.Range("F8").Value = 1.078362073794092
$ws.Range("I8").Value = 1.04326194892602
$ws.Range("J8").Value = 1.078955657950978
$ws.Range("K8").Value = 1.06172516073609
$ws.Range("L8").Value = 1.077637330559785
$ws.Range("M8").Value = 1.080980237282447
$ws.Range("N8").Value = 1.080487898743212
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.068886694976517
$ws.Range("D9").Value = 1.056779237326097
$ws.Range("E9").Value = 1.070671894287559
$ws.Range("F9").Value = 1.074718488429433
$ws.Range("I9").Value = 1.042591398176745
$ws.Range("J9").Value = 1.074487352854433
$ws.Range("K9").Value = 1.059875487909906
$ws.Range("L9").Value = 1.073725366898407
$ws.Range("M9").Value = 1.077759727608713
$ws.Range("N9").Value = 1.076013248140908
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.065345232014051
$ws.Range("D10").Value = 1.055245931631867
$ws.Range("E10").Value = 1.067755506243974
$ws.Range("F10").Value = 1.072267396266602
$ws.Range("I10").Value = 1.042131036652171
$ws.Range("J10").Value = 1.071481269185644
$ws.Range("K10").Value = 1.058624147248563
$ws.Range("L10").Value = 1.071091355166349
$ws.Range("M10").Value = 1.075588230818389
$ws.Range("N10").Value = 1.073002895488507
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.063804892556802
$ws.Range("D11").Value = 1.054578282783592
$ws.Range("E11").Value = 1.06648666445966
$ws.Range("F11").Value = 1.071200562347597
$ws.Range("I11").Value = 1.041928477889593
$ws.Range("J11").Value = 1.070172783664886
$ws.Range("K11").Value = 1.058077852880615
$ws.Range("L11").Value = 1.069944310246823
$ws.Range("M11").Value = 1.074641895967232
$ws.Range("N11").Value = 1.071692551768225
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.063231673440057
$ws.Range("D12").Value = 1.054329718373124
$ws.Range("E12").Value = 1.066014425830442
$ws.Range("F12").Value = 1.07080344480243
$ws.Range("I12").Value = 1.041852749981141
$ws.Range("J12").Value = 1.069685696003625
$ws.Range("K12").Value = 1.057874254303075
$ws.Range("L12").Value = 1.069517242862396
$ws.Range("M12").Value = 1.074289453033465
$ws.Range("N12").Value = 1.071204772386614
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.063354679843549
$ws.Range("D13").Value = 1.054383062257455
$ws.Range("E13").Value = 1.066115765366761
$ws.Range("F13").Value = 1.070888666597659
$ws.Range("I13").Value = 1.041869016063055
$ws.Range("J13").Value = 1.069790226343113
$ws.Range("K13").Value = 1.057917957864869
$ws.Range("L13").Value = 1.069608896132835
$ws.Range("M13").Value = 1.074365095760073
$ws.Range("N13").Value = 1.071309451171174
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.063757532100696
$ws.Range("D14").Value = 1.054557748064426
$ws.Range("E14").Value = 1.066447648295325
$ws.Range("F14").Value = 1.071167753915721
$ws.Range("I14").Value = 1.041922228192409
$ws.Range("J14").Value = 1.070132542611479
$ws.Range("K14").Value = 1.058061037306232
$ws.Range("L14").Value = 1.069909029360391
$ws.Range("M14").Value = 1.074612782061647
$ws.Range("N14").Value = 1.071652253567906
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.064005599961054
$ws.Range("D15").Value = 1.054665301893227
$ws.Range("E15").Value = 1.066652007646046
$ws.Range("F15").Value = 1.071339595788809
$ws.Range("I15").Value = 1.041954949053244
$ws.Range("J15").Value = 1.070343313831195
$ws.Range("K15").Value = 1.058149102760003
$ws.Range("L15").Value = 1.070093817508942
$ws.Range("M15").Value = 1.074765265622147
$ws.Range("N15").Value = 1.071863324106929
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.06544731150211
$ws.Range("D16").Value = 1.055290161997227
$ws.Range("E16").Value = 1.067839585668003
$ws.Range("F16").Value = 1.07233808084752
$ws.Range("I16").Value = 1.04214441156773
$ws.Range("J16").Value = 1.07156796244952
$ws.Range("K16").Value = 1.058660308187229
$ws.Range("L16").Value = 1.07116734149022
$ws.Range("M16").Value = 1.075650906413862
$ws.Range("N16").Value = 1.073089711866761
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.06634979550647
$ws.Range("D17").Value = 1.055681116996293
$ws.Range("E17").Value = 1.068582888271954
$ws.Range("F17").Value = 1.072962917328242
$ws.Range("I17").Value = 1.042262391085002
$ws.Range("J17").Value = 1.072334301096649
$ws.Range("K17").Value = 1.058979773162923
$ws.Range("L17").Value = 1.071838975291798
$ws.Range("M17").Value = 1.076204807527998
$ws.Range("N17").Value = 1.07385713880267
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.066875539570612
$ws.Range("D18").Value = 1.05590879654239
$ws.Range("E18").Value = 1.069015864607074
$ws.Range("F18").Value = 1.07332684523449
$ws.Range("I18").Value = 1.042330896273649
$ws.Range("J18").Value = 1.072780636482088
$ws.Range("K18").Value = 1.059165682538415
$ws.Range("L18").Value = 1.072230102584525
$ws.Range("M18").Value = 1.076527305192947
$ws.Range("N18").Value = 1.074304108035551
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.067054693956595
$ws.Range("D19").Value = 1.055986369021443
$ws.Range("E19").Value = 1.069163400922431
$ws.Range("F19").Value = 1.073450846251036
$ws.Range("I19").Value = 1.042354202293307
$ws.Range("J19").Value = 1.072932714825361
$ws.Range("K19").Value = 1.059229000375065
$ws.Range("L19").Value = 1.072363361632017
$ws.Range("M19").Value = 1.076637170230733
$ws.Range("N19").Value = 1.074456402347516
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.066253036007212
$ws.Range("D20").Value = 1.055639208325255
$ws.Range("E20").Value = 1.068503199085087
$ws.Range("F20").Value = 1.072895933128249
$ws.Range("I20").Value = 1.042249765127192
$ws.Range("J20").Value = 1.072252148393929
$ws.Range("K20").Value = 1.058945542062404
$ws.Range("L20").Value = 1.071766980218837
$ws.Range("M20").Value = 1.076145439685701
$ws.Range("N20").Value = 1.073774869433689
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.063638931949987
$ws.Range("D21").Value = 1.054506323272885
$ws.Range("E21").Value = 1.066349943070145
$ws.Range("F21").Value = 1.071085593245152
$ws.Range("I21").Value = 1.041906572073314
$ws.Range("J21").Value = 1.070031768460766
$ws.Range("K21").Value = 1.058018922853011
$ws.Range("L21").Value = 1.069820675536752
$ws.Range("M21").Value = 1.074539870520133
$ws.Range("N21").Value = 1.07155133630634
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.061989140230612
$ws.Range("D22").Value = 1.053790730405412
$ws.Range("E22").Value = 1.064990685366184
$ws.Range("F22").Value = 1.069942443509787
$ws.Range("I22").Value = 1.041687963791746
$ws.Range("J22").Value = 1.068629590926973
$ws.Range("K22").Value = 1.057432377618445
$ws.Range("L22").Value = 1.068591134893955
$ws.Range("M22").Value = 1.073524982230447
$ws.Range("N22").Value = 1.070147167519616
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.062864326446042
$ws.Range("D23").Value = 1.054170396678463
$ws.Range("E23").Value = 1.065711777125742
$ws.Range("F23").Value = 1.070548922477027
$ws.Range("I23").Value = 1.041804122017262
$ws.Range("J23").Value = 1.069373504098967
$ws.Range("K23").Value = 1.057743693855322
$ws.Range("L23").Value = 1.069243498674253
$ws.Range("M23").Value = 1.074063513055716
$ws.Range("N23").Value = 1.070892137133643
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.066296759496704
$ws.Range("D24").Value = 1.055658146155433
$ws.Range("E24").Value = 1.068539208988122
$ws.Range("F24").Value = 1.072926202037761
$ws.Range("I24").Value = 1.042255471212656
$ws.Range("J24").Value = 1.072289271693013
$ws.Range("K24").Value = 1.05896101094904
$ws.Range("L24").Value = 1.071799513624677
$ws.Range("M24").Value = 1.076172267260477
$ws.Range("N24").Value = 1.073812045452116
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.070254160674468
$ws.Range("D25").Value = 1.057370596020469
$ws.Range("E25").Value = 1.07179765241344
$ws.Range("F25").Value = 1.075664241495245
$ws.Range("I25").Value = 1.042767082332178
$ws.Range("J25").Value = 1.075647193343127
$ws.Range("K25").Value = 1.060356843178307
$ws.Range("L25").Value = 1.074741188789528
$ws.Range("M25").Value = 1.078596542653539
$ws.Range("N25").Value = 1.077174735736132
